$d = $word.ActiveDocument
$s = $d.Styles.Add("Temp1", 1)
$s.Delete()
